$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 577; this pushes existing rows 577-649 down to 578-650
$ws.Rows("577").Insert()

# Populate the new row 577 with the new record's data
$ws.Range("A577").Value = 3
$ws.Range("B577").Value = "Femacal de La Calera"
$ws.Range("C577").Value = "Coquimbo"
$ws.Range("D577").Value = 45142
$ws.Range("E577").Value = 5
$ws.Range("F577").Value = 100112031
$ws.Range("G577").Value = "Poroto verde"
$ws.Range("H577").Value = "Sin especificar"
$ws.Range("I577").Value = "Primera"
$ws.Range("J577").Value = 40
$ws.Range("K577").Value = 35000
$ws.Range("L577").Value = 35000
$ws.Range("M577").Value = 35000
$ws.Range("N577").Value = "`$/malla 25 kilos"
$ws.Range("O577").Value = "Perú"
$ws.Range("P577").Value = 1400
$ws.Range("Q577").Value = 25
$ws.Range("R577").Value = "Hortaliza"
